$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "NumOptions" (col M) and "Model_Base" (col N) columns were swapped
# for the header row and the populated data rows (2-7).
for ($r = 2; $r -le 7; $r++) {
    $mCell = $ws.Cells.Item($r, 13)   # column M
    $nCell = $ws.Cells.Item($r, 14)   # column N

    $mVal = $mCell.Value()
    $nVal = $nCell.Value()

    $mCell.Value = $nVal
    $nCell.Value = $mVal
}

# Update the active selection to reflect where editing ended up.
$ws.Range("O4").Select()
